$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.602.58'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '3.444.63'
$ws.Range('E3').Value = '  -2.34%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'590.71"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.63%  '
$ws.Range('D6').Value = "'178.39"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.95%  '
$ws.Range('E7').Value = '  +1.60%  '
$ws.Range('D9').Value = '3.445.80'
$ws.Range('E9').Value = '  -2.26%  '
$ws.Range('E10').Value = '  -1.62%  '
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('E12').Value = '  -3.71%  '
$ws.Range('D13').Value = '4.041.47'
$ws.Range('E13').Value = '  -2.41%  '
$ws.Range('D14').Value = "'31.93"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.83%  '
$ws.Range('E15').Value = '  -1.08%  '
$ws.Range('D16').Value = '67.623.98'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D18').Value = '3.444.33'
$ws.Range('E18').Value = '  -2.50%  '
$ws.Range('E19').Value = '  -4.42%  '
$ws.Range('D20').Value = "'13.97"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.66%  '
$ws.Range('D21').Value = "'389.18"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('D22').Value = "'7.83"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.69%  '
$ws.Range('E23').Value = '  +2.00%  '
$ws.Range('D24').Value = "'0.999"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.20%  '
$ws.Range('E25').Value = '  -2.70%  '
$ws.Range('D26').Value = "'71.38"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.20%  '
$ws.Range('E27').Value = '  -5.39%  '
$ws.Range('D28').Value = "'10.22"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.16%  '
$ws.Range('E29').Value = '  -2.20%  '
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('E31').Value = '  -4.17%  '
$ws.Range('E32').Value = '  -2.14%  '
$ws.Range('E33').Value = '  -5.97%  '
$ws.Range('E34').Value = '  -3.69%  '
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('E36').Value = '  -3.77%  '
$ws.Range('E37').Value = '  -7.72%  '
$ws.Range('D38').Value = "'160.98"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.70%  '
$ws.Range('D39').Value = "'0.884"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.38%  '
$ws.Range('E40').Value = '  -5.77%  '
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('E42').Value = '  -3.88%  '
$ws.Range('D43').Value = "'6.63"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.29%  '
$ws.Range('D44').Value = "'25.80"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.00%  '
$ws.Range('D45').Value = "'0.0713"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.88%  '
$ws.Range('D46').Value = "'25.94"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.12%  '
$ws.Range('D47').Value = '2.693.72'
$ws.Range('E47').Value = '  -6.52%  '
$ws.Range('E48').Value = '  -3.10%  '
$ws.Range('E49').Value = '  -3.68%  '
$ws.Range('D50').Value = "'323.65"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -8.04%  '
$ws.Range('E51').Value = '  -5.06%  '
